# Auto-generated script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.973.79"
$ws.Range("E2").Value = "  -1.09%  "
$ws.Range("D3").Value = "3.174.23"
$ws.Range("E3").Value = "  -4.34%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.04%  "
$ws.Range("D8").Value = "3.168.11"
$ws.Range("E8").Value = "  -4.50%  "
$ws.Range("E9").Value = "  -1.45%  "
$ws.Range("E10").Value = "  -4.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.29"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.458"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000238"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.14%  "
$ws.Range("D15").Value = "3.693.29"
$ws.Range("E15").Value = "  -4.39%  "
$ws.Range("E16").Value = "  -1.99%  "
$ws.Range("D17").Value = "3.174.22"
$ws.Range("E17").Value = "  -4.25%  "
$ws.Range("D18").Value = "62.866.68"
$ws.Range("E18").Value = "  -1.44%  "
$ws.Range("E19").Value = "  -3.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "460.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.20%  "
$ws.Range("E21").Value = "  -1.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.712"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.91%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.79"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.77%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.81%  "
$ws.Range("E31").Value = "  -6.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.104"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.93%  "
$ws.Range("E34").Value = "  -6.37%  "
$ws.Range("E35").Value = "  -5.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.86"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.14"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.05%  "
$ws.Range("D38").Value = "0.0₃0708"
$ws.Range("E38").Value = "  -4.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0387"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.60%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "405.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.45%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.69"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.78%  "
$ws.Range("E43").Value = "  -3.69%  "
$ws.Range("D44").Value = "2.790.15"
$ws.Range("E44").Value = "  -9.63%  "
$ws.Range("E45").Value = "  -4.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.998"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.95%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.69"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.87%  "
$ws.Range("E51").Value = "  -2.54%  "
